# Applies the "Översikt LUDVIKA" update:
#   1. Column C ("Förändrad") for existing rows 2-408 moves from 45189 (2023-09-20)
#      to 45190 (2023-09-21).
#   2. Row 408 gains an explicit 15pt custom row height.
#   3. A brand new record is appended as row 409 (case "A 44520-2023").
#   4. The sheet dimension grows from A1:Y408 to A1:Y409 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the "Förändrad" date for every existing data row (2..408).
$ws.Range("C2:C408").Value2 = 45190

# 2. Row 408 picks up an explicit row height in the new file.
$ws.Rows.Item(408).RowHeight = 15

# 3. Append the new case as row 409.
$ws.Range("A409").Value2 = "A 44520-2023"

$ws.Range("B409").Value2 = 45189
$ws.Range("B409").NumberFormat = "YYYY-MM-DD"

$ws.Range("C409").Value2 = 45190
$ws.Range("C409").NumberFormat = "YYYY-MM-DD"

$ws.Range("D409").Value2 = "DALARNAS LÄN"
$ws.Range("E409").Value2 = "LUDVIKA"
# (no F409 - Markägare is blank for this case, matching row 408's pattern)

$ws.Range("G409").Value2 = 0.7
$ws.Range("H409:Q409").Value2 = 0

# R column carries the (empty) wrapped species-list text, same style as above rows.
$ws.Range("R409").Value2 = ""
$ws.Range("R409").WrapText = $true
